$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date label in A3
$ws.Range("A3").Value = "Date:28.05.19"

# Update sales data values
$ws.Range("B6").Value = 121160
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 8

$ws.Range("B7").Value = 72980
$ws.Range("C7").Value = 48
$ws.Range("D7").Value = 8

$ws.Range("B8").Value = 39345
$ws.Range("C8").Value = 33
$ws.Range("D8").Value = 8

# Update the selected cell in the sheet view
$ws.Range("E11").Select()
